$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: built from a format-copy of row 2 (A:AD), then values changed ---
# Row 2 already carries the exact style pattern that row 4 needs (A-I plus the
# J:AD "filler" cells), so copy formats down and then overwrite the cell
# values/clear the two cells that differ from row 2.
$ws.Range("A2:AD2").Copy()
$ws.Range("A4:AD4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A4").Value = "DeLuca Jr., William F"
$ws.Range("B4").Value = "Principal"
$ws.Range("C4").Value = "#6789"
$ws.Range("D4").Value = "MD"
$ws.Range("E4").Value = "0102/0304"
$ws.Range("F4").ClearContents()
$ws.Range("I4").ClearContents()

# --- Row 5 ---
$ws.Range("G2").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)  # xlPasteFormats -> style 2
$ws.Range("F2").Copy()
$ws.Range("E5:F5").PasteSpecial(-4122)  # xlPasteFormats -> style 5
$ws.Range("G2").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H5").PasteSpecial(-4122)     # xlPasteFormats -> style 3
$ws.Range("G2").Copy()
$ws.Range("I5").PasteSpecial(-4122)

$ws.Range("A5").Value = "brown, malcolm d "
$ws.Range("B5").Value = "Sub"
$ws.Range("C5").Value = "#11667"
$ws.Range("D5").Value = "MD"

# --- Row 6 ---
$ws.Range("G2").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)  # style 2
$ws.Range("F2").Copy()
$ws.Range("E6:F6").PasteSpecial(-4122)  # style 5
$ws.Range("G2").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("I6").PasteSpecial(-4122)

$ws.Range("A6").Value = "David C. Kauffman"
$ws.Range("B6").Value = "Sub"

# --- Row 7 ---
$ws.Range("G2").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)  # style 2
$ws.Range("F2").Copy()
$ws.Range("C7").PasteSpecial(-4122)     # style 5
$ws.Range("G2").Copy()
$ws.Range("D7").PasteSpecial(-4122)     # style 2
$ws.Range("F2").Copy()
$ws.Range("E7:F7").PasteSpecial(-4122)  # style 5
$ws.Range("G2").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("I7").PasteSpecial(-4122)

$ws.Range("A7").Value = "Liang, Cheng Yi"
$ws.Range("B7").Value = "Sub"

$ws.Range("A7").Select() | Out-Null

$wb.Windows.Item(1).Left = 930
$wb.Windows.Item(1).Top = 120
